$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column CN (92), shifting "nom"/"url_produit" right.
$ws.Columns("CN").Insert()

# New header timestamp for the freshly inserted column
$ws.Range("CN1").Value = "2026-01-31 20:12:41"

# Carry forward last known price snapshot into the new column
for ($r = 2; $r -le 206; $r++) {
    $ws.Cells.Item($r, 92).Value = $ws.Cells.Item($r, 91).Value2
}
